# Bulk invoicing sales main/merge template - format & header label changes
# - B1 header: "Billing ID" -> "Settlement ID"
# - C1 header: "Actual Billing ID" -> "Billing ID"
# - Sample row (row 2) font color for columns A & B switches from the
#   dimmed gray (#444444) used for the rest of the font set to the
#   standard automatic/theme text color (matches the black already used
#   by column C's sample + the placeholder rows below).
# - Selection moves to B1 (the newly relabeled settlement id header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text changes -------------------------------------------------
$ws.Range("B1").Value = "Settlement ID"
$ws.Range("C1").Value = "Billing ID"

# --- Font color changes ---------------------------------------------------
# Row 2's A2/B2 sample cells were a dim gray; make them automatic/theme text
# color (black) like the rest of the sample row (C2) and placeholder rows.
$ws.Range("A2:B2").Font.ThemeColor = 1

# --- Selection -------------------------------------------------------------
[void]$ws.Range("B1").Select()
